# Add example values to mandatory ENA templates
# - bump Version on the isa_template sheet from 1.0.0 -> 1.0.1
# - fill in example values for the mandatory columns of row 2 on the
#   "New Table" data sheet (B,C,D,H,I,J,K,L,M,N,O,P,Q,R,S)

$wb = $excel.ActiveWorkbook

# --- isa_template sheet: bump Version value (row 4, col B) ---
$wsMeta = $wb.Worksheets.Item("isa_template")
$wsMeta.Range("B4").Value = "1.0.1"

# --- New Table sheet: populate example values on row 2 ---
$wsData = $wb.Worksheets.Item("New Table")

$wsData.Range("B2").Value = "Illumina HiSeq 1500"
$wsData.Range("C2").Value = "GENEPIO"
$wsData.Range("D2").Value = "http://purl.obolibrary.org/obo/GENEPIO_0100115"

$wsData.Range("H2").Value = "genomic DNA"
$wsData.Range("I2").Value = "EFO"
$wsData.Range("J2").Value = "http://purl.obolibrary.org/obo/EFO_0008479"

$wsData.Range("K2").Value = "Polymerase Chain Reaction"
$wsData.Range("L2").Value = "NCIT"
$wsData.Range("M2").Value = "http://purl.obolibrary.org/obo/NCIT_C17003"

$wsData.Range("N2").Value = "Whole Genome Sequencing"
$wsData.Range("O2").Value = "NCIT"
$wsData.Range("P2").Value = "http://purl.obolibrary.org/obo/NCIT_C101294"

$wsData.Range("Q2").Value = "single-end"
$wsData.Range("R2").Value = "DPBO"
$wsData.Range("S2").Value = "http://purl.obolibrary.org/obo/DPBO_0000086"
